# Updates the NATMI ligand-receptor pair sheet (Spp1-Itga4) so that the target-cluster
# cross join covers all four clusters (ECs, FAPs, M2, sCs) for every sending cluster,
# instead of only M2/sCs, per "Natmi following Dr Hou advice". This replaces the existing
# 8 data rows (rows 2-9) and appends 8 new rows (rows 10-17), for 16 rows total (A2:T17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

# Row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Spp1"
$data[0,2] = "Itga4"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 463.732605
$data[0,7] = 1391.197815
$data[0,8] = 0.3632113435366598
$data[0,9] = 0.3632113435366598
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 23.65990166666667
$data[0,13] = 70.979705
$data[0,14] = 0.2997993941754699
$data[0,15] = 0.29979939417547
$data[0,16] = 10971.86783392717
$data[0,17] = 98746.81050534456
$data[0,18] = 0.1088905407499491
$data[0,19] = 0.1088905407499491

# Row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Spp1"
$data[1,2] = "Itga4"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 463.732605
$data[1,7] = 1391.197815
$data[1,8] = 0.3632113435366598
$data[1,9] = 0.3632113435366598
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.07690566666666666
$data[1,13] = 0.230717
$data[1,14] = 0.0009744872400636476
$data[1,15] = 0.0009744872400636479
$data[1,16] = 35.66366514259499
$data[1,17] = 320.972986283355
$data[1,18] = 0.000353944819722849
$data[1,19] = 0.0003539448197228491

# Row 4: ECs -> M2
$data[2,0] = "ECs"
$data[2,1] = "Spp1"
$data[2,2] = "Itga4"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 463.732605
$data[2,7] = 1391.197815
$data[2,8] = 0.3632113435366598
$data[2,9] = 0.3632113435366598
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 53.21452433333334
$data[2,13] = 159.643573
$data[2,14] = 0.6742919890890982
$data[2,15] = 0.6742919890890983
$data[2,16] = 24677.30999293255
$data[2,17] = 222095.789936393
$data[2,18] = 0.2449104992930581
$data[2,19] = 0.2449104992930582

# Row 5: ECs -> sCs
$data[3,0] = "ECs"
$data[3,1] = "Spp1"
$data[3,2] = "Itga4"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 463.732605
$data[3,7] = 1391.197815
$data[3,8] = 0.3632113435366598
$data[3,9] = 0.3632113435366598
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 1.967779333333333
$data[3,13] = 5.903338
$data[3,14] = 0.02493412949536815
$data[3,15] = 0.02493412949536816
$data[3,16] = 912.52343631183
$data[3,17] = 8212.710926806469
$data[3,18] = 0.009056358673929723
$data[3,19] = 0.009056358673929727

# Row 6: FAPs -> ECs
$data[4,0] = "FAPs"
$data[4,1] = "Spp1"
$data[4,2] = "Itga4"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 12.24662533333333
$data[4,7] = 36.739876
$data[4,8] = 0.009591978638444229
$data[4,9] = 0.009591978638444227
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 23.65990166666667
$data[4,13] = 70.979705
$data[4,14] = 0.2997993941754699
$data[4,15] = 0.29979939417547
$data[4,16] = 289.7539511351756
$data[4,17] = 2607.78556021658
$data[4,18] = 0.002875669384749629
$data[4,19] = 0.002875669384749629

# Row 7: FAPs -> FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Spp1"
$data[5,2] = "Itga4"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 12.24662533333333
$data[5,7] = 36.739876
$data[5,8] = 0.009591978638444229
$data[5,9] = 0.009591978638444227
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.07690566666666666
$data[5,13] = 0.230717
$data[5,14] = 0.0009744872400636476
$data[5,15] = 0.0009744872400636479
$data[5,16] = 0.9418348856768889
$data[5,17] = 8.476513971092
$data[5,18] = 0.00000934726079012698
$data[5,19] = 0.000009347260790126982

# Row 8: FAPs -> M2
$data[6,0] = "FAPs"
$data[6,1] = "Spp1"
$data[6,2] = "Itga4"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 12.24662533333333
$data[6,7] = 36.739876
$data[6,8] = 0.009591978638444229
$data[6,9] = 0.009591978638444227
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 53.21452433333334
$data[6,13] = 159.643573
$data[6,14] = 0.6742919890890982
$data[6,15] = 0.6742919890890983
$data[6,16] = 651.6983418018832
$data[6,17] = 5865.285076216948
$data[6,18] = 0.006467794355416698
$data[6,19] = 0.006467794355416698

# Row 9: FAPs -> sCs
$data[7,0] = "FAPs"
$data[7,1] = "Spp1"
$data[7,2] = "Itga4"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 12.24662533333333
$data[7,7] = 36.739876
$data[7,8] = 0.009591978638444229
$data[7,9] = 0.009591978638444227
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.967779333333333
$data[7,13] = 5.903338
$data[7,14] = 0.02493412949536815
$data[7,15] = 0.02493412949536816
$data[7,16] = 24.09865623400978
$data[7,17] = 216.887906106088
$data[7,18] = 0.0002391676374877735
$data[7,19] = 0.0002391676374877735

# Row 10: M2 -> ECs
$data[8,0] = "M2"
$data[8,1] = "Spp1"
$data[8,2] = "Itga4"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 689.7685036666667
$data[8,7] = 2069.305511
$data[8,8] = 0.5402504422695089
$data[8,9] = 0.5402504422695089
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 23.65990166666667
$data[8,13] = 70.979705
$data[8,14] = 0.2997993941754699
$data[8,15] = 0.29979939417547
$data[8,16] = 16319.85496951714
$data[8,17] = 146878.6947256542
$data[8,18] = 0.1619667552954284
$data[8,19] = 0.1619667552954285

# Row 11: M2 -> FAPs
$data[9,0] = "M2"
$data[9,1] = "Spp1"
$data[9,2] = "Itga4"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 689.7685036666667
$data[9,7] = 2069.305511
$data[9,8] = 0.5402504422695089
$data[9,9] = 0.5402504422695089
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.07690566666666666
$data[9,13] = 0.230717
$data[9,14] = 0.0009744872400636476
$data[9,15] = 0.0009744872400636479
$data[9,16] = 53.04710662015411
$data[9,17] = 477.423959581387
$data[9,18] = 0.0005264671624303787
$data[9,19] = 0.0005264671624303789

# Row 12: M2 -> M2
$data[10,0] = "M2"
$data[10,1] = "Spp1"
$data[10,2] = "Itga4"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 689.7685036666667
$data[10,7] = 2069.305511
$data[10,8] = 0.5402504422695089
$data[10,9] = 0.5402504422695089
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 53.21452433333334
$data[10,13] = 159.643573
$data[10,14] = 0.6742919890890982
$data[10,15] = 0.6742919890890983
$data[10,16] = 36705.70282273676
$data[10,17] = 330351.3254046308
$data[10,18] = 0.3642865453241722
$data[10,19] = 0.3642865453241722

# Row 13: M2 -> sCs
$data[11,0] = "M2"
$data[11,1] = "Spp1"
$data[11,2] = "Itga4"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 689.7685036666667
$data[11,7] = 2069.305511
$data[11,8] = 0.5402504422695089
$data[11,9] = 0.5402504422695089
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 1.967779333333333
$data[11,13] = 5.903338
$data[11,14] = 0.02493412949536815
$data[11,15] = 0.02493412949536816
$data[11,16] = 1357.312206299524
$data[11,17] = 12215.80985669572
$data[11,18] = 0.01347067448747785
$data[11,19] = 0.01347067448747785

# Row 14: sCs -> ECs
$data[12,0] = "sCs"
$data[12,1] = "Spp1"
$data[12,2] = "Itga4"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 111.00921
$data[12,7] = 333.02763
$data[12,8] = 0.08694623555538696
$data[12,9] = 0.08694623555538696
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 23.65990166666667
$data[12,13] = 70.979705
$data[12,14] = 0.2997993941754699
$data[12,15] = 0.29979939417547
$data[12,16] = 2626.46699269435
$data[12,17] = 23638.20293424915
$data[12,18] = 0.02606642874534271
$data[12,19] = 0.02606642874534272

# Row 15: sCs -> FAPs
$data[13,0] = "sCs"
$data[13,1] = "Spp1"
$data[13,2] = "Itga4"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 111.00921
$data[13,7] = 333.02763
$data[13,8] = 0.08694623555538696
$data[13,9] = 0.08694623555538696
$data[13,10] = 1
$data[13,11] = 0.3333333333333333
$data[13,12] = 0.07690566666666666
$data[13,13] = 0.230717
$data[13,14] = 0.0009744872400636476
$data[13,15] = 0.0009744872400636479
$data[13,16] = 8.53723730119
$data[13,17] = 76.83513571071
$data[13,18] = 0.00008472799712029282
$data[13,19] = 0.00008472799712029285

# Row 16: sCs -> M2
$data[14,0] = "sCs"
$data[14,1] = "Spp1"
$data[14,2] = "Itga4"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 111.00921
$data[14,7] = 333.02763
$data[14,8] = 0.08694623555538696
$data[14,9] = 0.08694623555538696
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 53.21452433333334
$data[14,13] = 159.643573
$data[14,14] = 0.6742919890890982
$data[14,15] = 0.6742919890890983
$data[14,16] = 5907.30230676911
$data[14,17] = 53165.72076092199
$data[14,18] = 0.05862715011645114
$data[14,19] = 0.05862715011645115

# Row 17: sCs -> sCs
$data[15,0] = "sCs"
$data[15,1] = "Spp1"
$data[15,2] = "Itga4"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 111.00921
$data[15,7] = 333.02763
$data[15,8] = 0.08694623555538696
$data[15,9] = 0.08694623555538696
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 1.967779333333333
$data[15,13] = 5.903338
$data[15,14] = 0.02493412949536815
$data[15,15] = 0.02493412949536816
$data[15,16] = 218.44162924766
$data[15,17] = 1965.97466322894
$data[15,18] = 0.002167928696472801
$data[15,19] = 0.002167928696472801

$ws.Range("A2:T17").Value = $data
